# clv_inputs_v1.1.xlsx maintenance edit:
#  - Inputs sheet: the old "Description" column (C) is relabelled "Key" and
#    now holds the snake_case identifier for each parameter; a brand new
#    column D ("Description") takes over the old human-readable text. Two
#    new "Monthly Churn" rows are inserted before CAC Target (which becomes
#    row 8).
#  - Workbook-level defined names: replace the old Title-Case pair with a
#    lowercase, snake_case set (one per Inputs parameter).
#  - Add a README sheet with short documentation notes.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inputs")

# ---------------------------------------------------------------------
# 1) Inputs sheet: add a new "Description" header in D1, styled like the
#    existing header cells, then relabel C1 from "Description" to "Key".
# ---------------------------------------------------------------------
$ws.Range("D1").Value = "Description"
$ws.Range("C1").Copy()
$ws.Range("D1").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Range("C1").Value = "Key"

# ---------------------------------------------------------------------
# 2) Move the old human-readable text out of column C and into the new
#    column D, then put the snake_case key into column C. Also rename the
#    "Basic Price"/"Pro Price" parameter labels.
# ---------------------------------------------------------------------
$ws.Range("D2").Value = "Margin after COGS"
$ws.Range("C2").Value = "gross_margin"

$ws.Range("D3").Value = "Annual WACC"
$ws.Range("C3").Value = "discount_rate"

$ws.Range("A4").Value = "Price Basic"
$ws.Range("D4").Value = "Monthly Basic Price"
$ws.Range("C4").Value = "price_basic"

$ws.Range("A5").Value = "Price Pro"
$ws.Range("D5").Value = "Monthly Pro Price"
$ws.Range("C5").Value = "price_pro"

# ---------------------------------------------------------------------
# 3) Insert the two new "Monthly Churn" rows (6 & 7), pushing CAC Target
#    down to row 8.
# ---------------------------------------------------------------------
$ws.Range("A8").Value = "CAC Target"
$ws.Range("B8").Value = 150
$ws.Range("C8").Value = "cac_target"
$ws.Range("D8").Value = "Target Cost per Acq"

$ws.Range("A6").Value = "Monthly Churn Basic"
$ws.Range("B6").Value = 0.05
$ws.Range("C6").Value = "monthly_churn_basic"
$ws.Range("D6").Value = "Est. Monthly Churn"

$ws.Range("A7").Value = "Monthly Churn Pro"
$ws.Range("B7").Value = 0.02
$ws.Range("C7").Value = "monthly_churn_pro"
$ws.Range("D7").Value = "Est. Monthly Churn"

# ---------------------------------------------------------------------
# 4) Workbook defined names: drop the old Title-Case pair, add the new
#    lowercase snake_case set.
# ---------------------------------------------------------------------
$oldNames = @()
foreach ($n in $wb.Names) {
    $oldNames += $n.Name
}
foreach ($nm in $oldNames) {
    $wb.Names.Item($nm).Delete()
}

$wb.Names.Add("gross_margin", '=Inputs!$B$2')
$wb.Names.Add("discount_rate", '=Inputs!$B$3')
$wb.Names.Add("price_basic", '=Inputs!$B$4')
$wb.Names.Add("price_pro", '=Inputs!$B$5')
$wb.Names.Add("monthly_churn_basic", '=Inputs!$B$6')
$wb.Names.Add("monthly_churn_pro", '=Inputs!$B$7')
$wb.Names.Add("cac_target", '=Inputs!$B$8')

# ---------------------------------------------------------------------
# 5) Add a README sheet (after the last existing sheet) with short
#    documentation notes.
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$readme = $wb.Worksheets.Add($null, $lastSheet)
$readme.Name = "README"
$readme.Range("A1").Value = "This file contains inputs for the CLV model."
$readme.Range("A2").Value = "Do not rename sheets."
$readme.Range("A3").Value = "Named Ranges defined: gross_margin, discount_rate, price_basic, price_pro, etc."

# Restore "Inputs" as the selected/active sheet (it was tabSelected before
# the edit; adding README would otherwise steal that flag).
$ws.Activate()

Write-Output "done"
